# Adds "green hydrogen" and "low carbon hydrogen" rows to the BFoHfC sheet,
# mirroring the existing "hydrogen" row (row 11) which drives every year
# column (B:AK, 2015-2050) off the Data sheet's capacity-weighted CHP
# fraction (Data!$A$6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BFoHfC")

# --- New row 12: "green hydrogen" ---
$ws.Range("A12").Value = "green hydrogen"
for ($col = 2; $col -le 37; $col++) {
    $ws.Cells.Item(12, $col).Formula = '=Data!$A$6'
}
# Clear the "inherited from precedent" number format so these cells stay
# General, matching the unstyled cells in the rest of the fuel-type block.
$ws.Range("B12:AK12").Style = "Normal"

# --- New row 13: "low carbon hydrogen" ---
$ws.Range("A13").Value = "low carbon hydrogen"
for ($col = 2; $col -le 37; $col++) {
    $ws.Cells.Item(13, $col).Formula = '=Data!$A$6'
}
$ws.Range("B13:AK13").Style = "Normal"

# The edit was made directly on the BFoHfC sheet, leaving the selection on
# the row right after the new data (A14), and that sheet/cell becomes the
# one active when the workbook is saved.
$ws.Activate() | Out-Null
$ws.Range("A14").Select() | Out-Null
